$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.937.53"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "2.524.20"
$ws.Range("E3").Value = "  -2.18%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.26"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.55"
$ws.Range("E6").Value = "  +3.91%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("D9").Value = "2.523.55"
$ws.Range("E9").Value = "  -2.15%  "
$ws.Range("E10").Value = "  +1.17%  "
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.13"
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.342"
$ws.Range("E13").Value = "  -4.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.54"
$ws.Range("E14").Value = "  -0.76%  "
$ws.Range("D15").Value = "2.983.41"
$ws.Range("E15").Value = "  -2.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000176"
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("D17").Value = "66.851.92"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "2.530.68"
$ws.Range("E18").Value = "  -1.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.02"
$ws.Range("E19").Value = "  +3.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.32"
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "354.30"
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.17"
$ws.Range("E22").Value = "  -1.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.60"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  +5.09%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.76"
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.00"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "2.650.33"
$ws.Range("E29").Value = "  -2.27%  "
$ws.Range("E30").Value = "  -1.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "532.11"
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.09"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.84"
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.46"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.63"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.57"
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.44"
$ws.Range("E40").Value = "  +1.12%  "
$ws.Range("E41").Value = "  -1.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.79"
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.11"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("E45").Value = "  +2.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "149.05"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.555"
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("D48").Value = "0.0₆0276"
$ws.Range("E48").Value = "  -3.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.67"
$ws.Range("E49").Value = "  -1.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.68"
$ws.Range("E50").Value = "  -1.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0755"
$ws.Range("E51").Value = "  -0.61%  "
